# Generate Report for Handback
# Update the timestamp strings recorded on the handback-status workbook's
# "Overview", "zh-cn" and "de-de" sheets to reflect the latest report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for
# 3d8830c2-7274-4dd3-8e97-dffc123f6e14.md
$wsOverview.Range("G3").Value = "2016-11-09 05:49:49"

# zh-cn sheet: Handoff / Handback datetimes for
# 3d8830c2-7274-4dd3-8e97-dffc123f6e14.7f104da06406dc0a2315aacab9a4bd525f6e9604.zh-cn.xlf
$wsZhCn.Range("H3").Value = "2016-11-09 05:49:36"
$wsZhCn.Range("K3").Value = "2016-11-09 05:50:34"

# de-de sheet: Handback datetime for
# 3d8830c2-7274-4dd3-8e97-dffc123f6e14.7f104da06406dc0a2315aacab9a4bd525f6e9604.de-de.xlf
$wsDeDe.Range("K3").Value = "2016-11-09 05:50:54"
